$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# TP du 18 mars 2024 FSIL - add new journal entry row (row 26)
$ws.Range("A26").Value = "18/03/2024"
$ws.Range("B26").Value = "FSIL"
$ws.Range("C26").Value = "TP"
$ws.Range("E26").Value = "x"
$ws.Range("G26").Value = "Fin question score pour tous (fix#4). 4 ont commencé en autonomie scrabble score."
